# Update cryptocurrency price (D) and volume-change (E) figures for the latest run.
# Values that look numeric are stored as text (matching the source data, which uses
# a dotted thousands style e.g. "20.408.91"), so a leading apostrophe forces Excel to
# keep them as text instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '20.408.91'
$ws.Range("E2").Value = '  +2.36%  '
$ws.Range("D3").Value = '1.461.99'
$ws.Range("E3").Value = '  +3.42%  '
$ws.Range("E4").Value = '  +0.80%  '
$ws.Range("D5").Value = '''0.9505'
$ws.Range("E5").Value = '  -5.01%  '
$ws.Range("E6").Value = '  -0.45%  '
$ws.Range("D7").Value = '''0.3652'
$ws.Range("E7").Value = '  -0.30%  '
$ws.Range("D8").Value = '''0.3069'
$ws.Range("E8").Value = '  -0.71%  '
$ws.Range("D9").Value = '''39.81'
$ws.Range("E9").Value = '  +0.18%  '
$ws.Range("D10").Value = '''1.035'
$ws.Range("E10").Value = '  +0.55%  '
$ws.Range("D11").Value = '''0.06582'
$ws.Range("E11").Value = '  +0.96%  '
$ws.Range("E12").Value = '  -0.06%  '
$ws.Range("D13").Value = '''5.429'
$ws.Range("E13").Value = '  -1.02%  '
$ws.Range("D14").Value = '''17.95'
$ws.Range("E14").Value = '  +2.29%  '
$ws.Range("D15").Value = '''6.140'
$ws.Range("E15").Value = '  -0.71%  '
$ws.Range("D16").Value = '''0.00001024'
$ws.Range("E16").Value = '  +0.69%  '
$ws.Range("D17").Value = '1.461.63'
$ws.Range("E17").Value = '  +3.53%  '
$ws.Range("D18").Value = '''0.9703'
$ws.Range("E18").Value = '  -3.04%  '
$ws.Range("D19").Value = '''0.05813'
$ws.Range("E19").Value = '  +2.67%  '
$ws.Range("D20").Value = '''69.41'
$ws.Range("E20").Value = '  -2.38%  '
$ws.Range("D21").Value = '''5.439'
$ws.Range("E21").Value = '  -3.04%  '
$ws.Range("E22").Value = '  -1.58%  '
$ws.Range("D23").Value = '''10.90'
$ws.Range("E23").Value = '  +0.22%  '
$ws.Range("D24").Value = '''2.247'
$ws.Range("E24").Value = '  +0.45%  '
$ws.Range("D25").Value = '20.436.41'
$ws.Range("E25").Value = '  +2.41%  '
$ws.Range("D26").Value = '''141.72'
$ws.Range("E26").Value = '  +6.68%  '
$ws.Range("D27").Value = '''2.084'
$ws.Range("E27").Value = '  -7.77%  '
$ws.Range("D28").Value = '''17.12'
$ws.Range("E28").Value = '  -0.69%  '
$ws.Range("D29").Value = '1.614.02'
$ws.Range("E29").Value = '  +2.73%  '
$ws.Range("D30").Value = '''112.26'
$ws.Range("E30").Value = '  +2.39%  '
$ws.Range("D31").Value = '''3.842'
$ws.Range("E31").Value = '  -0.78%  '
$ws.Range("D32").Value = '''4.889'
$ws.Range("E32").Value = '  -6.86%  '
$ws.Range("D33").Value = '''0.07895'
$ws.Range("E33").Value = '  +2.80%  '
$ws.Range("D34").Value = '''0.7900'
$ws.Range("E34").Value = '  -3.07%  '
$ws.Range("E35").Value = '  +3.17%  '
$ws.Range("D36").Value = '''0.05726'
$ws.Range("E36").Value = '  -0.72%  '
$ws.Range("D37").Value = '''1.149'
$ws.Range("E37").Value = '  +4.87%  '
$ws.Range("D38").Value = '''4.682'
$ws.Range("E38").Value = '  -4.68%  '
$ws.Range("D39").Value = '''0.02026'
$ws.Range("E39").Value = '  -0.99%  '
$ws.Range("D40").Value = '''0.9581'
$ws.Range("E40").Value = '  -4.02%  '
$ws.Range("D41").Value = '''10.33'
$ws.Range("E41").Value = '  -0.82%  '
$ws.Range("D42").Value = '''7.459'
$ws.Range("E42").Value = '  -10.15%  '
$ws.Range("D43").Value = '''0.1858'
$ws.Range("E43").Value = '  -1.30%  '
$ws.Range("D44").Value = '''0.5260'
$ws.Range("E44").Value = '  -0.76%  '
$ws.Range("E45").Value = '  -1.32%  '
$ws.Range("D46").Value = '''11.92'
$ws.Range("E46").Value = '  -3.31%  '
$ws.Range("D47").Value = '''117.06'
$ws.Range("E47").Value = '  +1.31%  '
$ws.Range("D48").Value = '''0.5145'
$ws.Range("E48").Value = '  -0.57%  '
$ws.Range("E49").Value = '  -1.17%  '
$ws.Range("D50").Value = '''0.06418'
$ws.Range("E50").Value = '  +3.83%  '
$ws.Range("D51").Value = '''0.9907'
$ws.Range("E51").Value = '  -1.03%  '
